# Scheduled-runner price/profit refresh for Siren_Profits workbook.
# Updates currentAveragePrice* (H:K) and Leve cost/profit (L:N) columns
# for the affected Leve rows on each crafting-job sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 31: Hush Little Wailer / Weak Silencing Potion
$ws.Range("H31").Value = 53.5
$ws.Range("I31").Value = 58.22222
$ws.Range("K31").Value = 174.66666
$ws.Range("M31").Value = 55.33333999999999

# ALC row 58: A Matter of Vital Importance / Mega-Potion of Vitality
$ws.Range("H58").Value = 537.4286
$ws.Range("I58").Value = 330.8421
$ws.Range("K58").Value = 992.5263
$ws.Range("M58").Value = -842.5263

# ALC row 75: Tomes Roam on the Range / Dhalmelskin Codex
$ws.Range("H75").Value = 196665.67
$ws.Range("J75").Value = 196665.67
$ws.Range("L75").Value = 196665.67
$ws.Range("N75").Value = -198537.67

# ALC row 78: Field Trip to the Unknown (L) / Dhalmelskin Codex
$ws.Range("H78").Value = 196665.67
$ws.Range("J78").Value = 196665.67
$ws.Range("L78").Value = 589997.01
$ws.Range("N78").Value = -599357.01

# ALC row 96: Scroll Down / Grade 1 Reisui of Intelligence
$ws.Range("H96").Value = 1062.1666
$ws.Range("I96").Value = 678.4286
$ws.Range("J96").Value = 1599.4
$ws.Range("K96").Value = 2035.2858
$ws.Range("L96").Value = 4798.200000000001
$ws.Range("M96").Value = -662.2857999999999
$ws.Range("N96").Value = -7544.200000000001

# ALC row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 33745.5
$ws.Range("I98").Value = 43039.92
$ws.Range("J98").Value = 20320.223
$ws.Range("K98").Value = 43039.92
$ws.Range("L98").Value = 20320.223
$ws.Range("M98").Value = -41541.92
$ws.Range("N98").Value = -23316.223

# ALC row 100: Asking for a Friend / Beetle Glue
$ws.Range("H100").Value = 53508.605
$ws.Range("I100").Value = 53276.05
$ws.Range("J100").Value = 53999.555
$ws.Range("K100").Value = 53276.05
$ws.Range("L100").Value = 53999.555
$ws.Range("M100").Value = -52735.05
$ws.Range("N100").Value = -55081.555

# ALC row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 33745.5
$ws.Range("I122").Value = 43039.92
$ws.Range("J122").Value = 20320.223
$ws.Range("K122").Value = 129119.76
$ws.Range("L122").Value = 60960.66900000001
$ws.Range("M122").Value = -126669.76
$ws.Range("N122").Value = -65860.66900000001

# ALC row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 3048.2195
$ws.Range("I132").Value = 3040.8108
$ws.Range("J132").Value = 3116.75
$ws.Range("K132").Value = 9122.432400000002
$ws.Range("L132").Value = 9350.25
$ws.Range("M132").Value = -6592.432400000002
$ws.Range("N132").Value = -14410.25

# ALC row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 8484.839
$ws.Range("I137").Value = 9229.259
$ws.Range("J137").Value = 3460
$ws.Range("K137").Value = 27687.777
$ws.Range("L137").Value = 10380
$ws.Range("M137").Value = -25137.777
$ws.Range("N137").Value = -15480

$ws = $wb.Worksheets.Item("BSM")
# BSM row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 3292.7
$ws.Range("I86").Value = 3748.348
$ws.Range("J86").Value = 1795.5714
$ws.Range("K86").Value = 3748.348
$ws.Range("L86").Value = 1795.5714
$ws.Range("M86").Value = -2625.348
$ws.Range("N86").Value = -4041.5714

# BSM row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 3292.7
$ws.Range("I89").Value = 3748.348
$ws.Range("J89").Value = 1795.5714
$ws.Range("K89").Value = 18741.74
$ws.Range("L89").Value = 8977.857
$ws.Range("M89").Value = -13125.74
$ws.Range("N89").Value = -20209.857

# BSM row 131: Plying with Precision / Chondrite Pliers
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# BSM row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 2266.087
$ws.Range("I134").Value = 1463.55
$ws.Range("J134").Value = 7616.3335
$ws.Range("K134").Value = 4390.65
$ws.Range("L134").Value = 22849.0005
$ws.Range("M134").Value = -1855.65
$ws.Range("N134").Value = -27919.0005

$ws = $wb.Worksheets.Item("CRP")
# CRP row 29: Grinding It Out / Mudstone Grinding Wheel
$ws.Range("H29").Value = 1000
$ws.Range("J29").Value = 1000
$ws.Range("L29").Value = 1000
$ws.Range("N29").Value = -1586

# CRP row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 2522.6206
$ws.Range("I31").Value = 1906.5
$ws.Range("J31").Value = 5480
$ws.Range("K31").Value = 1906.5
$ws.Range("L31").Value = 5480
$ws.Range("M31").Value = -1611.5
$ws.Range("N31").Value = -6070

# CRP row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 2522.6206
$ws.Range("I34").Value = 1906.5
$ws.Range("J34").Value = 5480
$ws.Range("K34").Value = 1906.5
$ws.Range("L34").Value = 5480
$ws.Range("M34").Value = -1704.5
$ws.Range("N34").Value = -5884

# CRP row 86: Birch, Please / Birch Lumber
$ws.Range("H86").Value = 15197.8
$ws.Range("I86").Value = 13329.667
$ws.Range("K86").Value = 13329.667
$ws.Range("M86").Value = -12206.667

# CRP row 89: Built This City on Blocks and Soul (L) / Birch Lumber
$ws.Range("H89").Value = 15197.8
$ws.Range("I89").Value = 13329.667
$ws.Range("K89").Value = 66648.33499999999
$ws.Range("M89").Value = -61032.33499999999

# CRP row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 21114.285
$ws.Range("I122").Value = 26960
$ws.Range("K122").Value = 80880
$ws.Range("M122").Value = -78430

# CRP row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 3352.7827
$ws.Range("I134").Value = 1758.5834
$ws.Range("K134").Value = 5275.7502
$ws.Range("M134").Value = -2740.7502

$ws = $wb.Worksheets.Item("CUL")
# CUL row 4: In Hot Water / Boiled Egg
$ws.Range("H4").Value = 19549188
$ws.Range("J4").Value = 420696960
$ws.Range("L4").Value = 1262090880
$ws.Range("N4").Value = -1262091104

# CUL row 63: The Next to Last Supper / Stuffed Cabbage Rolls
$ws.Range("H63").Value = 2498
$ws.Range("J63").Value = 2997
$ws.Range("L63").Value = 8991
$ws.Range("N63").Value = -10489

# CUL row 66: Nostalgia through the Stomach (L) / Stuffed Cabbage Rolls
$ws.Range("H66").Value = 2498
$ws.Range("J66").Value = 2997
$ws.Range("L66").Value = 26973
$ws.Range("N66").Value = -34461

# CUL row 92: Oh No Udon / Gyr Abanian Flour
$ws.Range("H92").Value = 307
$ws.Range("I92").Value = 268.4
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 805.1999999999999
$ws.Range("L92").Value = 1500
$ws.Range("M92").Value = 442.8000000000001
$ws.Range("N92").Value = -3996

$ws = $wb.Worksheets.Item("GSM")
# GSM row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 20913.5
$ws.Range("I122").Value = 19200
$ws.Range("J122").Value = 22627
$ws.Range("K122").Value = 57600
$ws.Range("L122").Value = 67881
$ws.Range("M122").Value = -55150
$ws.Range("N122").Value = -72781

# GSM row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 33897.9
$ws.Range("J126").Value = 24427.143
$ws.Range("L126").Value = 73281.429
$ws.Range("N126").Value = -78221.429

# GSM row 128: To Fight at Her Side / Manganese Rapier
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# GSM row 131: Star Athletes / Star Quartz Wristband of Aiming
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# GSM row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 3039.6943
$ws.Range("I132").Value = 2710.3438
$ws.Range("J132").Value = 5674.5
$ws.Range("K132").Value = 8131.0314
$ws.Range("L132").Value = 17023.5
$ws.Range("M132").Value = -5601.0314
$ws.Range("N132").Value = -22083.5

$ws = $wb.Worksheets.Item("LTW")
# LTW row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 29116.941
$ws.Range("I7").Value = 41272
$ws.Range("K7").Value = 41272
$ws.Range("M7").Value = -41160

# LTW row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 29116.941
$ws.Range("I126").Value = 41272
$ws.Range("K126").Value = 123816
$ws.Range("M126").Value = -121346

$ws = $wb.Worksheets.Item("WVR")
# WVR row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 5373.643
$ws.Range("I122").Value = 4073.15
$ws.Range("K122").Value = 12219.45
$ws.Range("M122").Value = -9769.450000000001

# WVR row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 10488.733
$ws.Range("I132").Value = 11107.784
$ws.Range("K132").Value = 33323.352
$ws.Range("M132").Value = -30793.352

# WVR row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 1589.3636
$ws.Range("I136").Value = 1060.625
$ws.Range("J136").Value = 2999.3333
$ws.Range("K136").Value = 3181.875
$ws.Range("L136").Value = 8997.999899999999
$ws.Range("M136").Value = -631.875
$ws.Range("N136").Value = -14097.9999
